# Applies the "NIT-9015967185" statement-of-account update:
#  - adds a second overdue period (2506) as a new table row
#  - relabels the existing row's period as 2507 (most recent month)
#  - updates the totals (VALOR MORA, Cant. Periodos) accordingly
#  - swaps the "Novedad de Ingreso" / "Novedad de Retiro" header order

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new data row right after the existing worker row (16) ---
# Copying the whole row first makes Excel clone formatting/borders into the
# freshly inserted row, same as a user would do via "Copy row > Insert copied cells".
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert()
$excel.CutCopyMode = 0

# --- 2. Fix up the two period cells: old row keeps "2507", new row gets "2506" ---
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"

# --- 3. Update the summary totals ---
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2

# --- 4. Swap the "Novedad de Ingreso" / "Novedad de Retiro" column headers ---
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

Write-Host ("dim=" + $ws.UsedRange.Address())
Write-Host ("E16=" + $ws.Range("E16").Text)
Write-Host ("E17=" + $ws.Range("E17").Text)
Write-Host ("E11=" + $ws.Range("E11").Text)
Write-Host ("F13=" + $ws.Range("F13").Text)
Write-Host ("H15=" + $ws.Range("H15").Text)
Write-Host ("I15=" + $ws.Range("I15").Text)
